$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '68.356.79'
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +2.06%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '3.573.30'
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +0.38%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.12%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '206.79'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +8.29%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '561.65'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -1.20%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.607'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -0.97%  '
$ws.Cells.Item(8, 5).Value = '  +0.14%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.672'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.21%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '63.23'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +13.64%  '
$ws.Cells.Item(11, 5).Value = '  -1.98%  '
$ws.Cells.Item(12, 5).Value = '  +3.47%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '10.08'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +2.67%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '4.144.14'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +0.39%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '3.575.11'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.27%  '
$ws.Cells.Item(16, 5).Value = '  +0.10%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '19.07'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +5.21%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '68.167.65'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +1.89%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '12.12'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.16%  '
$ws.Cells.Item(20, 5).Value = '  +0.28%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '396.94'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -0.67%  '
$ws.Cells.Item(22, 2).Value = 'PancakeSwap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '4.15'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.25%  '
$ws.Cells.Item(23, 2).Value = 'RenderToken'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '12.37'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +4.61%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '84.02'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -1.74%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '2.86'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -1.16%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '12.37'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.45%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '3.85'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +4.72%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '9.07'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.27%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '713.11'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +10.95%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '31.30'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +0.70%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '7.51'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -3.44%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '12.01'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.61%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '63.57'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.20%  '
$ws.Cells.Item(34, 5).Value = '  -1.53%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '41.16'
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -2.26%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.421'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +4.44%  '
$ws.Cells.Item(37, 5).Value = '  +0.00%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '3.20'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +5.71%  '
$ws.Cells.Item(39, 5).Value = '  +29.12%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '3.150.87'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -0.65%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0728'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -4.01%  '
$ws.Cells.Item(42, 5).Value = '  -1.21%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -0.16%  '
$ws.Cells.Item(44, 5).Value = '  -3.70%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +8.94%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '0.0410'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.36%  '
$ws.Cells.Item(47, 5).Value = '  +0.02%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '3.06'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.91%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '8.62'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +1.56%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '138.69'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -2.41%  '
$ws.Cells.Item(51, 5).Value = '  +0.28%  '
